$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "60.930.24"
Set-TextValue "E2" "  +4.08%  "
Set-TextValue "D3" "2.713.95"
Set-TextValue "E3" "  +3.76%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "529.39"
Set-TextValue "E5" "  +1.27%  "
Set-TextValue "D6" "147.50"
Set-TextValue "E6" "  +2.25%  "
Set-TextValue "E7" "  -0.29%  "
Set-TextValue "D8" "0.581"
Set-TextValue "E8" "  +2.45%  "
Set-TextValue "D9" "2.739.41"
Set-TextValue "E9" "  +4.83%  "
Set-TextValue "D10" "7.23"
Set-TextValue "E10" "  +15.52%  "
Set-TextValue "E11" "  +2.31%  "
Set-TextValue "D12" "0.342"
Set-TextValue "E12" "  +2.95%  "
Set-TextValue "E13" "  +3.34%  "
Set-TextValue "D14" "3.188.09"
Set-TextValue "E14" "  +3.79%  "
Set-TextValue "D15" "60.939.91"
Set-TextValue "E15" "  +4.07%  "
Set-TextValue "B16" "WrappedEther"
Set-TextValue "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "2.838.51"
Set-TextValue "E16" "  +8.22%  "
Set-TextValue "B17" "Avalanche"
Set-TextValue "C17" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D17" "21.53"
Set-TextValue "E17" "  +4.08%  "
Set-TextValue "D18" "0.0000139"
Set-TextValue "E18" "  +2.35%  "
Set-TextValue "B19" "BitcoinCash"
Set-TextValue "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "347.03"
Set-TextValue "E19" "  +0.45%  "
Set-TextValue "B20" "Polkadot"
Set-TextValue "C20" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D20" "4.54"
Set-TextValue "E20" "  +2.80%  "
Set-TextValue "D21" "10.62"
Set-TextValue "E21" "  +4.78%  "
Set-TextValue "D22" "6.46"
Set-TextValue "E22" "  +5.71%  "
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  +0.17%  "
Set-TextValue "E24" "  +3.81%  "
Set-TextValue "B25" "Kaspa"
Set-TextValue "C25" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D25" "0.172"
Set-TextValue "E25" "  +5.38%  "
Set-TextValue "B26" "Polygon"
Set-TextValue "C26" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D26" "0.420"
Set-TextValue "E26" "  +1.92%  "
Set-TextValue "D27" "0.992"
Set-TextValue "E27" "  -0.61%  "
Set-TextValue "D28" "0.0₃0828"
Set-TextValue "E28" "  +4.18%  "
Set-TextValue "D29" "7.36"
Set-TextValue "E29" "  +5.69%  "
Set-TextValue "D30" "6.80"
Set-TextValue "E30" "  +10.31%  "
Set-TextValue "E31" "  -0.23%  "
Set-TextValue "E32" "  +2.21%  "
Set-TextValue "D33" "19.11"
Set-TextValue "E33" "  +1.79%  "
Set-TextValue "D34" "149.97"
Set-TextValue "E34" "  +0.12%  "
Set-TextValue "E35" "  +8.26%  "
Set-TextValue "D36" "1.25"
Set-TextValue "E36" "  +10.35%  "
Set-TextValue "D37" "0.927"
Set-TextValue "E37" "  -4.67%  "
Set-TextValue "D38" "0.898"
Set-TextValue "E38" "  +8.04%  "
Set-TextValue "D39" "1.54"
Set-TextValue "E39" "  +9.11%  "
Set-TextValue "D40" "37.27"
Set-TextValue "E40" "  +1.94%  "
Set-TextValue "E41" "  +2.63%  "
Set-TextValue "D42" "284.20"
Set-TextValue "E42" "  +1.75%  "
Set-TextValue "D43" "20.36"
Set-TextValue "E43" "  +4.86%  "
Set-TextValue "B44" "Stellar"
Set-TextValue "C44" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D44" "0.0993"
Set-TextValue "E44" "  +1.59%  "
Set-TextValue "B45" "Mantle"
Set-TextValue "C45" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D45" "0.614"
Set-TextValue "E45" "  +3.41%  "
Set-TextValue "E46" "  -0.07%  "
Set-TextValue "D47" "2.131.63"
Set-TextValue "E47" "  +7.74%  "
Set-TextValue "D48" "5.01"
Set-TextValue "E48" "  +8.79%  "
Set-TextValue "D49" "0.0548"
Set-TextValue "E49" "  +5.59%  "
Set-TextValue "D50" "10.54"
Set-TextValue "E50" "  +2.23%  "
Set-TextValue "D51" "19.44"
Set-TextValue "E51" "  +6.63%  "
